# Append four new paragraphs after the "Histórico de coleta de lixo nas ruas -
# Confiabilidade" paragraph and before the document's final (empty) paragraph:
#   1. an empty paragraph
#   2. "Questão 5"
#   3. the answer text
#   4. an empty paragraph
#
# We do this by repeatedly inserting a new paragraph immediately before the
# last paragraph of the document (which duplicates its formatting, matching
# the style used throughout the rest of the document), then filling in the
# text for the two paragraphs that need it.

$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range

# Insert 4 new empty paragraphs right before the existing last paragraph.
$r.InsertParagraphBefore() | Out-Null
$r.InsertParagraphBefore() | Out-Null
$r.InsertParagraphBefore() | Out-Null
$r.InsertParagraphBefore() | Out-Null

$total = $d.Paragraphs.Count

# Paragraph indices (1-based) of the four newly inserted paragraphs, in
# document order: total-4 .. total-1. The original (untouched) last
# paragraph is now at index $total.
$pEmpty1 = $d.Paragraphs.Item($total - 4)
$pQuestion = $d.Paragraphs.Item($total - 3)
$pAnswer = $d.Paragraphs.Item($total - 2)
$pEmpty2 = $d.Paragraphs.Item($total - 1)

$pQuestion.Range.Font.NameAscii = "Calibri"
$pQuestion.Range.Font.NameFarEast = "Calibri"
$pQuestion.Range.Font.NameOther = "Calibri"
$pQuestion.Range.Font.NameBi = "Calibri"
$pQuestion.Range.Text = "Questão 5"

$pAnswer.Range.Font.NameAscii = "Calibri"
$pAnswer.Range.Font.NameFarEast = "Calibri"
$pAnswer.Range.Font.NameOther = "Calibri"
$pAnswer.Range.Font.NameBi = "Calibri"
$pAnswer.Range.Text = "Como estamos definindo qual método de projeto, quais são as bases de dados necessárias, como a equipe será organizada, é um período de análise, logo é planejamento de escopo e abordagem."
